$d = $word.ActiveDocument

# Locate the paragraph that currently holds the "retirement-beach" URL
# (the last paragraph of the body, right before the sectPr) and drop
# the insertion point at its very end.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$r = $lastPara.Range
$r.Collapse(0)

# Insert a blank paragraph, then a paragraph carrying the new citation URL.
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.MoveStart(1, 1)

$r.InsertParagraphAfter()
$r.Collapse(0)
$r.MoveStart(1, 1)

$r.InsertAfter("https://pixabay.com/photos/fish-sea-bream-barbecue-grilled-2366925/")

# Register the latent "Hyperlink" / "Unresolved Mention" character styles
# that Word stamps into styles.xml once a URL-like run has been typed.
$d.Styles.Add("Hyperlink", 2) | Out-Null
$d.Styles.Add("Unresolved Mention", 2) | Out-Null
